$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $newValue
    $r.Style = "Normal"
}

Set-TextValue 'D2' '25.855.28'
Set-TextValue 'E2' '  +0.22%  '
Set-TextValue 'D3' '1.741.37'
Set-TextValue 'E3' '  -0.44%  '
Set-TextValue 'E4' '  +0.05%  '
Set-TextValue 'D5' '231.04'
Set-TextValue 'E5' '  -2.13%  '
Set-TextValue 'E6' '  +0.07%  '
Set-TextValue 'D7' '0.5166'
Set-TextValue 'E7' '  +1.55%  '
Set-TextValue 'D8' '0.2794'
Set-TextValue 'E8' '  +4.64%  '
Set-TextValue 'D9' '39.45'
Set-TextValue 'E9' '  -3.13%  '
Set-TextValue 'D10' '0.06092'
Set-TextValue 'E10' '  -1.61%  '
Set-TextValue 'D11' '1.754.48'
Set-TextValue 'E11' '  +0.27%  '
Set-TextValue 'D12' '0.07039'
Set-TextValue 'E12' '  +1.49%  '
Set-TextValue 'E13' '  -1.27%  '
Set-TextValue 'D14' '0.6381'
Set-TextValue 'E14' '  +1.68%  '
Set-TextValue 'D15' '4.501'
Set-TextValue 'E15' '  +0.60%  '
Set-TextValue 'D16' '76.93'
Set-TextValue 'E16' '  -0.97%  '
Set-TextValue 'D17' '1.001'
Set-TextValue 'E17' '  +0.10%  '
Set-TextValue 'E18' '  +0.03%  '
Set-TextValue 'D19' '25.865.33'
Set-TextValue 'E19' '  +0.20%  '
Set-TextValue 'E20' '  -1.81%  '
Set-TextValue 'D21' '0.000006562'
Set-TextValue 'E21' '  -1.53%  '
Set-TextValue 'D22' '1.974.65'
Set-TextValue 'E22' '  -0.11%  '
Set-TextValue 'D23' '4.124'
Set-TextValue 'E23' '  +1.66%  '
Set-TextValue 'D24' '8.611'
Set-TextValue 'E24' '  +4.27%  '
Set-TextValue 'D25' '5.130'
Set-TextValue 'E25' '  -0.02%  '
Set-TextValue 'D26' '139.50'
Set-TextValue 'E26' '  +2.03%  '
Set-TextValue 'D27' '1.511'
Set-TextValue 'E27' '  +3.85%  '
Set-TextValue 'D28' '15.04'
Set-TextValue 'E28' '  -0.41%  '
Set-TextValue 'D29' '1.806'
Set-TextValue 'D30' '101.99'
Set-TextValue 'E30' '  -0.48%  '
Set-TextValue 'D31' '0.08228'
Set-TextValue 'E31' '  +0.51%  '
Set-TextValue 'D32' '3.651'
Set-TextValue 'E32' '  -1.36%  '
Set-TextValue 'D33' '3.418'
Set-TextValue 'E33' '  +0.53%  '
Set-TextValue 'D34' '0.04481'
Set-TextValue 'E34' '  +1.43%  '
Set-TextValue 'D35' '2.614'
Set-TextValue 'E35' '  -1.40%  '
Set-TextValue 'D36' '0.9781'
Set-TextValue 'E36' '  -1.86%  '
Set-TextValue 'D37' '0.6129'
Set-TextValue 'E37' '  +2.34%  '
Set-TextValue 'E38' '  +0.00%  '
Set-TextValue 'D39' '0.01584'
Set-TextValue 'E39' '  +1.23%  '
Set-TextValue 'D40' '1.917'
Set-TextValue 'E40' '  -1.26%  '
Set-TextValue 'D41' '1.001'
Set-TextValue 'E41' '  -0.05%  '
Set-TextValue 'D42' '100.60'
Set-TextValue 'E42' '  -0.87%  '
Set-TextValue 'D43' '0.3824'
Set-TextValue 'E43' '  +0.08%  '
Set-TextValue 'B44' 'FraxShare'
Set-TextValue 'C44' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D44' '4.973'
Set-TextValue 'E44' '  +1.65%  '
Set-TextValue 'B45' 'TrustWalletToken'
Set-TextValue 'C45' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D45' '0.7204'
Set-TextValue 'E45' '  -4.13%  '
Set-TextValue 'D46' '0.05414'
Set-TextValue 'E46' '  -1.59%  '
Set-TextValue 'D47' '6.260'
Set-TextValue 'E47' '  +5.50%  '
Set-TextValue 'D48' '0.1118'
Set-TextValue 'E48' '  +1.91%  '
Set-TextValue 'D49' '53.13'
Set-TextValue 'E49' '  +0.84%  '
Set-TextValue 'D50' '7.627'
Set-TextValue 'D51' '29.74'
Set-TextValue 'E51' '  -1.04%  '
